# Update the "想去人数" (number of people interested) figures for two
# events on both the "展览" and "全部类型" worksheets, which contain
# duplicated data:
#   F3: 286 -> 287
#   F4: 19  -> 20
#   F5: 81  -> 84

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 287
    $ws.Range("F4").Value = 20
    $ws.Range("F5").Value = 84
}
